# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Los periodos de mora (columna E) se reordenan de forma ascendente
$ws.Range("E16").Value = "2504"
$ws.Range("E17").Value = "2505"
$ws.Range("E18").Value = "2506"
$ws.Range("E19").Value = "2507"

# Nuevo salario base actualizado para cada periodo (columna G)
$ws.Range("G16").Value = 1423500
$ws.Range("G17").Value = 1423500
$ws.Range("G18").Value = 1423500
$ws.Range("G19").Value = 1423500
